# Issu no.108 - Kolom M,N (SKP3) di takeout
#
# The "TANGGAL MULAI SKPPP" (M) and "TANGGAL SELESAI SKPPP" (N) columns are
# removed from the "PEND. PANDU" sheet. The remaining "NO SK" / "TANGGAL
# MULAI SK" header columns (K/L) are relabeled to "NO SK/PKL" / "TANGGAL
# MULAI SK/PKL" to reflect that they now also cover PKL documents.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove columns M:N (TANGGAL MULAI SKPPP / TANGGAL SELESAI SKPPP) entirely -
# this shifts every column from O onward two places to the left.
$ws.Columns("M:N").Delete()

# Relabel the remaining SK columns to cover SK/PKL.
$ws.Range("K1").Value = "NO SK/PKL"
$ws.Range("L1").Value = "TANGGAL MULAI SK/PKL"

# Match the author's final selection position.
$ws.Range("L2").Select()
